# The authored change removes the "Feuil1" worksheet from the workbook,
# leaving "Feuil2" ("IDENTIFICATION / Liste des étapes") as the sole sheet.
# (The workbook was re-uploaded after the sheet was deleted in Excel, which
# is why sheet count, shared-string usage counts, the active-tab pointer,
# etc. all shift down accordingly - that bookkeeping is handled
# automatically by Excel/the engine when the sheet is removed.)

$wb = $excel.ActiveWorkbook

# Suppress the "data may exist on the deleted sheet" confirmation prompt,
# exactly as a human / macro deleting a sheet via the UI or VBA would.
$excel.DisplayAlerts = $false

$sheetToRemove = $wb.Worksheets.Item("Feuil1")
$sheetToRemove.Delete() | Out-Null

$excel.DisplayAlerts = $true

# "Feuil2" is now the only (and therefore active) sheet in the workbook.
$wb.Worksheets.Item("Feuil2").Activate()

Write-Output ("Worksheets remaining: " + $wb.Worksheets.Count)
foreach ($sheet in $wb.Worksheets) {
    Write-Output (" - " + $sheet.Name)
}
